$d = $word.ActiveDocument

# The readme talks about "Project Scarlett" (the codename for Xbox Series
# X|S). Per the November GDK release notes, that codename reference is
# replaced with the shipped product name "Xbox Series X|S devkit", which
# splits the sentence

#   "If using Project Scarlett, set the active solution platform to ..."

# into three runs:

#   "If using " + "Xbox Series X|S devkit" + ", set the active solution platform to ..."

$target = $d.Content
$found = $target.Find.Execute("Project Scarlett", $false, $false, $false, $false, `
                               $false, $true, 1, $false, "", 0)

if ($found) {
    $target.Text = "Xbox Series X|S devkit"

    # Briefly toggling a run-level formatting property and then clearing it
    # again forces Word to materialize the replaced text as its own <w:r>
    # run, distinct from the "If using " run before it and the ", set the
    # active solution platform to " run after it - matching how the
    # sentence is now split across three separate runs.
    $target.Bold = 1
    $target.Bold = 0
}
